$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Metadata element"
$ws.Range("B1").Value = "Description"
